$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the longer descriptions that are being added
$ws.Columns.Item(3).ColumnWidth = 32.1

# D2: add Parameter-Daten "-" for the /register GET row, wrap text
$ws.Range("D2").Value = "-"
$ws.Range("D2").WrapText = $true

# Row 3: POST /register
$ws.Range("A3").Value = "POST"
$ws.Range("B3").Value = "/register"
$ws.Range("C3").Value = "Erstellt ein Benutzer"
$ws.Range("D3").Value = "JSON"

# Row 4: DELETE /privacy
$ws.Range("A4").Value = "DELETE"
$ws.Range("B4").Value = "/privacy"
$ws.Range("C4").Value = "Löscht das Konto des Benutzers"
$ws.Range("D4").Value = "JSON"

# Row 5: GET /genres
$ws.Range("A5").Value = "GET"
$ws.Range("B5").Value = "/genres"
$ws.Range("C5").Value = "Zeigt alle Genres von Bücher"
$ws.Range("D5").Value = "-"
$ws.Range("D5").WrapText = $true

# Row 6: POST /forums
$ws.Range("A6").Value = "POST"
$ws.Range("B6").Value = "/forums"
$ws.Range("C6").Value = "Erstellt ein Diskussion in Forums"
$ws.Range("D6").Value = "JSON"

# Row 7: DELETE /forums
$ws.Range("A7").Value = "DELETE"
$ws.Range("B7").Value = "/forums"
$ws.Range("C7").Value = "Löscht das Diskussion in Forums"
$ws.Range("D7").Value = "JSON"

# Row 8: PUT /languages
$ws.Range("A8").Value = "PUT"
$ws.Range("B8").Value = "/languages"
$ws.Range("C8").Value = "Updates die bestehende Sprache"
$ws.Range("D8").Value = "JSON"

# Row 9: GET /searchbar
$ws.Range("A9").Value = "GET"
$ws.Range("B9").Value = "/searchbar"
$ws.Range("C9").Value = "Such das Buch im Suchleiste"
$ws.Range("D9").Value = "-"

# Row 10: GET /rating
$ws.Range("A10").Value = "GET"
$ws.Range("B10").Value = "/rating"
$ws.Range("C10").Value = "Zeigt das Rating/Review für das Buch"
$ws.Range("D10").Value = "-"

# Row 11: POST /rating
$ws.Range("A11").Value = "POST"
$ws.Range("B11").Value = "/rating"
$ws.Range("C11").Value = "Erstellt ein Rating/Review für das Buch"
$ws.Range("D11").Value = "JSON"

# Row 12: DELETE /rating
$ws.Range("A12").Value = "DELETE"
$ws.Range("B12").Value = "/rating"
$ws.Range("C12").Value = "Löscht das Rating/Review für das Buch"
$ws.Range("D12").Value = "JSON"

# Row 13: GET /favourite
$ws.Range("A13").Value = "GET"
$ws.Range("B13").Value = "/favourite"
$ws.Range("C13").Value = "Zeigt alle deine Lieblingsbücher"
$ws.Range("D13").Value = "-"

# Row 14: POST /favourite
$ws.Range("A14").Value = "POST"
$ws.Range("B14").Value = "/favourite"
$ws.Range("C14").Value = "Speichert das Buch in deine Favourite"
$ws.Range("D14").Value = "JSON"

# Row 15: DELETE /favourite
$ws.Range("A15").Value = "DELETE"
$ws.Range("B15").Value = "/favourite"
$ws.Range("C15").Value = "Löscht das Favourite"
$ws.Range("D15").Value = "JSON"

# Row 16: GET /books
$ws.Range("A16").Value = "GET"
$ws.Range("B16").Value = "/books"
$ws.Range("C16").Value = "Zeigt alle Bücher"
$ws.Range("D16").Value = "-"

# Row 17: POST /books
$ws.Range("A17").Value = "POST"
$ws.Range("B17").Value = "/books"
$ws.Range("C17").Value = "Erstellt ein Buch "
$ws.Range("D17").Value = "JSON"

# Row 18: DELETE /books
$ws.Range("A18").Value = "DELETE"
$ws.Range("B18").Value = "/books"
$ws.Range("C18").Value = "Löscht das Biuch"
$ws.Range("D18").Value = "JSON"

# Selection / view state matching the final saved state
$ws.Range("D10,D14").Select()
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
